$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 32
$ws.Range("F3").Value = 8985
$ws.Range("F4").Value = 2700
$ws.Range("F5").Value = 962
$ws.Range("F6").Value = 320
$ws.Range("F7").Value = 835
$ws.Range("F8").Value = 712
$ws.Range("F9").Value = 129
$ws.Range("F12").Value = 902
$ws.Range("F13").Value = 3863
$ws.Range("F14").Value = 294
$ws.Range("F15").Value = 173
$ws.Range("F16").Value = 805
$ws.Range("F18").Value = 55
$ws.Range("F22").Value = 1398
$ws.Range("F23").Value = 1363
$ws.Range("F24").Value = 485
$ws.Range("F26").Value = 153
$ws.Range("F27").Value = 175
$ws.Range("F28").Value = 371
$ws.Range("F33").Value = 725
$ws.Range("F34").Value = 57
$ws.Range("F37").Value = 6
$ws.Range("F42").Value = 362
$ws.Range("F43").Value = 27
$ws.Range("F44").Value = 27

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 4

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 32
$ws.Range("F4").Value = 320
$ws.Range("F5").Value = 835
$ws.Range("F6").Value = 712
$ws.Range("F7").Value = 129
$ws.Range("F10").Value = 902
$ws.Range("F12").Value = 3863
$ws.Range("F13").Value = 294
$ws.Range("F14").Value = 173
$ws.Range("F16").Value = 4
$ws.Range("F17").Value = 805
$ws.Range("F21").Value = 55
$ws.Range("F26").Value = 1398
$ws.Range("F27").Value = 1363
$ws.Range("F28").Value = 485
$ws.Range("F30").Value = 153
$ws.Range("F31").Value = 175
$ws.Range("F33").Value = 371
$ws.Range("F37").Value = 725
$ws.Range("F38").Value = 57
$ws.Range("F41").Value = 6
$ws.Range("F45").Value = 362
$ws.Range("F46").Value = 27
$ws.Range("F47").Value = 27
